# Weekly refresh of the fruit/vegetable price sheet: the daily rows (2-29)
# get re-shuffled against the underlying weekly source order. Columns
# D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) move together as a
# per-row unit; the other columns (market/category/metadata) stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 29

# new row (index into this array, offset by $firstRow) <- source row to pull from
$sourceRow = @(3,28,23,21,17,8,20,11,4,22,15,14,16,10,29,27,19,26,6,24,25,12,5,2,9,7,13,18)

$cols = @("D","J","K","L","M","P")

# Phase 1: snapshot the current values for the columns that move, before any
# writes, so the permutation (which has several cycles) doesn't clobber a
# source row before it's been read.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Phase 2: write each row's new values from its mapped source row.
for ($i = 0; $i -lt ($lastRow - $firstRow + 1); $i++) {
    $destRow = $firstRow + $i
    $srcRow = $sourceRow[$i]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}
